# Group 1 Team CBA - MITR.xlsx
# The edit: clear the Year-0 "Software Development" cost (cell B15) on Sheet1,
# which cascades through the CBA's dependent SUM/NPV/IRR formulas, and update
# the sheet's saved view/selection state (scroll position + active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the Year 0 "Software Development" cost (was 100) -> blank.
$ws.Range("B15").ClearContents()

# Recalculate the workbook so all dependent formulas (SUM, NPV, IRR, etc.)
# pick up the new value immediately.
$excel.Calculate()

# Update the sheet view to match the saved cursor/scroll position.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("B16").Select()
